$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(79).Insert()

# Fill new row 79 - same as row 80 pattern but new D/J/K/L/M/P values
$ws.Range("A79").Value = 4
$ws.Range("B79").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C79").Value = "Los Lagos"
$ws.Range("D79").Value = 45209
$ws.Range("E79").Value = 10
$ws.Range("F79").Value = 100112012
$ws.Range("G79").Value = "Espinaca"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 35
$ws.Range("K79").Value = 13000
$ws.Range("L79").Value = 13000
$ws.Range("M79").Value = 13000
$ws.Range("N79").Value = "$/cuna 10 kilos"
$ws.Range("O79").Value = "Región Metropolitana"
$ws.Range("P79").Value = 1300
$ws.Range("Q79").Value = 10
$ws.Range("R79").Value = "Hortaliza"

Write-Host "Done"
